# edit.ps1 - reproduce the target diff via Word COM interop
#
# Summary of the change (see commit message / diff):
#  1. "Exceptions: ... agent ..." paragraph: the sentence that was split by
#     a <w:proofErr> pair around "page, but" is retyped as one run (no
#     visible text change).
#  2. "Entry condition:" (Supervisor checks potential intervention use
#     case) is reworded from "The supervisor must be on his personal
#     homepage." to "The supervisor must be logged in. Possibly he have
#     received a notification about a new available suggestion."
#  3. The following step ("The supervisor clicks on the Recent
#     Suggestions button.") is reworded to "... button or clicks on the
#     new suggestion notification." and the runs that made up the first
#     half of the sentence are retyped as one run.
#  4. The "_GoBack" bookmark moves from the end of the "no recent
#     suggestions" exception paragraph to the end of the step edited in
#     (3).
#
$d = $word.ActiveDocument

function Get-FoundRange($doc, $searchText) {
    $r = $doc.Content
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    return $r
}

# Force Word's run-merge/normalize pass over a Range by writing a throwaway
# placeholder value first and then writing the desired final text back in.
# (Setting .Text to the exact same text it already holds is treated as a
# no-op by this engine and will not coalesce runs / drop <w:proofErr>
# splits, so we always go through a placeholder first to guarantee a real
# mutation is recorded.)
function Set-RangeTextForced($doc, $range, $finalText) {
    $start = $range.Start
    $placeholder = "##PLACEHOLDER_" + [guid]::NewGuid().ToString("N") + "##"
    $range.Text = $placeholder
    $r2 = $doc.Range($start, $start + $placeholder.Length)
    $r2.Text = $finalText
    return $doc.Range($start, $start + $finalText.Length)
}

# ---------------------------------------------------------------------
# Edit 1: "Exceptions: ... agent ..." paragraph - collapse the runs that
# were split by <w:proofErr> around "page, but" into a single run with
# identical text (no visible text change, just a run/proofErr cleanup).
# ---------------------------------------------------------------------
$oldAgentText = " The username and password furnished by the agent are not correct. In this case, the system does not redirect the agent to his personal page, but notifies him that an error has been made and allows to input his username and password again."
$r1 = Get-FoundRange $d $oldAgentText
Set-RangeTextForced $d $r1 $oldAgentText | Out-Null

# ---------------------------------------------------------------------
# Edit 2: "Entry condition:" for "Supervisor checks potential
# intervention" - reword the sentence.
# ---------------------------------------------------------------------
$oldEntryCond = " The supervisor must be on his personal homepage."
$newEntryCond = " The supervisor must be logged in. Possibly he have received a notification about a new available suggestion. "
$r2 = Get-FoundRange $d $oldEntryCond
Set-RangeTextForced $d $r2 $newEntryCond | Out-Null

# ---------------------------------------------------------------------
# Edit 3: "The supervisor clicks on the 'Recent Suggestions' button."
# - merge the leading runs and extend the sentence.
# ---------------------------------------------------------------------
$oldLead = "The " + [char]0x2018 + [char]0x2019  # placeholder, unused
$oldLeadText = "The supervisor clicks on the " + [char]0x201C + "Recent "
$r3 = Get-FoundRange $d $oldLeadText
Set-RangeTextForced $d $r3 $oldLeadText | Out-Null

# Now update the tail: the run right after the underlined "Suggestions"
# run currently holds the text: <quote> button. - locate it unambiguously
# by anchoring off "Suggestions" (unique in the doc) and taking the range
# immediately following it rather than searching for the ambiguous
# "<quote> button." text (which also occurs elsewhere in the document).
$rSugg = Get-FoundRange $d "Suggestions"
$rSugg.Collapse(0) | Out-Null   # wdCollapseEnd
$quoteButtonLen = ([char]0x201D + " button.").Length
$rSugg.MoveEnd(1, $quoteButtonLen) | Out-Null  # wdCharacter
$r4 = $d.Range($rSugg.Start, $rSugg.End)
if ($r4.Text -ne ([char]0x201D + " button.")) {
    throw "Unexpected tail text: [$($r4.Text)]"
}
$newTail = [char]0x201D + " button or clicks on the new suggestion notification."
Set-RangeTextForced $d $r4 $newTail | Out-Null

Write-Output "done"
